# Update LTP (column B) and PREV (column C) values on the "ltp" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$values = @{
    2  = @(2519.3, 2509.5)
    3  = @(381.65, 384.3)
    4  = @(1541, 1519.55)
    5  = @(7410.65, 7363.2)
    6  = @(243.4, 238.5)
    7  = @(197.8, 196.25)
    8  = @(45319.65, 45045.15)
    9  = @(540, 530.75)
    10 = @(3370.9, 3402.45)
    11 = @(147.5, 145)
    12 = @(1262.05, 1254)
    13 = @(1448.85, 1439.6)
    14 = @(698.3, 703.5)
    15 = @(452.95, 451)
    16 = @(1567.55, 1565.75)
    17 = @(301.65, 299.9)
    18 = @(19872.45, 19768.7)
    19 = @(583.45, 579.05)
    20 = @(619.95, 617.35)
    21 = @(627.25, 614.9)
    22 = @(268.85, 263.7)
    23 = @(129.5, 130.15)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
